# Amelx-Lamp1.xlsx was regenerated from an updated TPM run: a new
# "Resolving-Mac" cluster now appears alongside ECs/FAPs/MuSCs, and the
# ligand/receptor expression-derived statistics were recomputed. Rewrite
# the data rows (A2:T13) with the refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2-13, columns A-T (new TPM-derived values).
$rows = @(
    ,("FAPs", "Amelx", "Lamp1", "ECs", 2, 0.6666666666666666, 0.6742976666666666, 2.022893, 0.3960292783506769, 0.3960292783506769, 3, 1, 59.97191866666666, 179.915756, 0.1189776078066642, 0.1189776078066642, 40.43892482245644, 363.9503234021079, 0.04711861615956309, 0.04711861615956309)
    ,("FAPs", "Amelx", "Lamp1", "FAPs", 2, 0.6666666666666666, 0.6742976666666666, 2.022893, 0.3960292783506769, 0.3960292783506769, 3, 1, 123.224251, 369.672753, 0.2444631910017033, 0.2444631910017033, 83.08982492604765, 747.808424334429, 0.09681458111570826, 0.09681458111570826)
    ,("FAPs", "Amelx", "Lamp1", "MuSCs", 2, 0.6666666666666666, 0.6742976666666666, 2.022893, 0.3960292783506769, 0.3960292783506769, 3, 1, 55.54748266666667, 166.642448, 0.1102000195140577, 0.1102000195140576, 37.45553795134044, 337.099841562064, 0.04364243420238276, 0.04364243420238276)
    ,("FAPs", "Amelx", "Lamp1", "Resolving-Mac", 2, 0.6666666666666666, 0.6742976666666666, 2.022893, 0.3960292783506769, 0.3960292783506769, 3, 1, 265.3168996666666, 795.950699, 0.5263591816775748, 0.5263591816775749, 178.9025663724674, 1610.123097352207, 0.2084536468730228, 0.2084536468730228)
    ,("MuSCs", "Amelx", "Lamp1", "ECs", 3, 1, 0.8625470000000001, 2.587641, 0.5065920925430184, 0.5065920925430184, 3, 1, 59.97191866666666, 179.915756, 0.1189776078066642, 0.1189776078066642, 51.72859853017733, 465.557386771596, 0.06027311530454058, 0.06027311530454058)
    ,("MuSCs", "Amelx", "Lamp1", "FAPs", 3, 1, 0.8625470000000001, 2.587641, 0.5065920925430184, 0.5065920925430184, 3, 1, 123.224251, 369.672753, 0.2444631910017033, 0.2444631910017033, 106.286708027297, 956.5803722456731, 0.1238431194792965, 0.1238431194792965)
    ,("MuSCs", "Amelx", "Lamp1", "MuSCs", 3, 1, 0.8625470000000001, 2.587641, 0.5065920925430184, 0.5065920925430184, 3, 1, 55.54748266666667, 166.642448, 0.1102000195140577, 0.1102000195140576, 47.91231453168534, 431.210830785168, 0.05582645848390793, 0.05582645848390792)
    ,("MuSCs", "Amelx", "Lamp1", "Resolving-Mac", 3, 1, 0.8625470000000001, 2.587641, 0.5065920925430184, 0.5065920925430184, 3, 1, 265.3168996666666, 795.950699, 0.5263591816775748, 0.5263591816775749, 228.8482958567843, 2059.634662711059, 0.2666493992752734, 0.2666493992752735)
    ,("Resolving-Mac", "Amelx", "Lamp1", "ECs", 1, 0.3333333333333333, 0.1658013333333333, 0.497404, 0.09737862910630474, 0.09737862910630474, 3, 1, 59.97191866666666, 179.915756, 0.1189776078066642, 0.1189776078066642, 9.943424077491555, 89.490816697424, 0.01158587634256054, 0.01158587634256054)
    ,("Resolving-Mac", "Amelx", "Lamp1", "FAPs", 1, 0.3333333333333333, 0.1658013333333333, 0.497404, 0.09737862910630474, 0.09737862910630474, 3, 1, 123.224251, 369.672753, 0.2444631910017033, 0.2444631910017033, 20.43074511480133, 183.876706033212, 0.0238054904066986, 0.0238054904066986)
    ,("Resolving-Mac", "Amelx", "Lamp1", "MuSCs", 1, 0.3333333333333333, 0.1658013333333333, 0.497404, 0.09737862910630474, 0.09737862910630474, 3, 1, 55.54748266666667, 166.642448, 0.1102000195140577, 0.1102000195140576, 9.209846689443555, 82.88862020499201, 0.01073112682776696, 0.01073112682776696)
    ,("Resolving-Mac", "Amelx", "Lamp1", "Resolving-Mac", 1, 0.3333333333333333, 0.1658013333333333, 0.497404, 0.09737862910630474, 0.09737862910630474, 3, 1, 265.3168996666666, 795.950699, 0.5263591816775748, 0.5263591816775749, 43.98989572059955, 395.909061485396, 0.05125613552927863, 0.05125613552927864)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
